$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Previously added")
$ws2 = $wb.Worksheets.Item("New")

# ---------------------------------------------------------------------------
# Step 1: remember the hyperlink targets currently on sheet "New" (rows 2-9,
# in row order), then drop all hyperlinks from that sheet. (Hyperlinks.Delete()
# clears every hyperlink on the worksheet it is invoked on.)
# ---------------------------------------------------------------------------
$oldLinks = @(
    "https://www.ss.com/msg/lv/real-estate/wood/cesis-and-reg/ligatnes-pag/bgnghf.html",
    "https://www.ss.com/msg/lv/real-estate/wood/daugavpils-and-reg/laucesas-pag/adhnd.html",
    "https://www.ss.com/msg/lv/real-estate/wood/jekabpils-and-reg/dunavas-pag/mbdni.html",
    "https://www.ss.com/msg/lv/real-estate/wood/kuldiga-and-reg/padures-pag/cghhpx.html",
    "https://www.ss.com/msg/lv/real-estate/wood/madona-and-reg/indranu-pag/bgcpkd.html",
    "https://www.ss.com/msg/lv/real-estate/wood/madona-and-reg/erglu-pag/ekgnc.html",
    "https://www.ss.com/msg/lv/real-estate/wood/preili-and-reg/turku-pag/cmcnb.html",
    "https://www.ss.com/msg/lv/real-estate/wood/rezekne-and-reg/ozolmuizas-pag/blmkl.html"
)

$ws2.Hyperlinks.Delete()

# ---------------------------------------------------------------------------
# Step 2: move the 8 existing data rows from "New" (rows 2-9) down to the
# bottom of "Previously added" (rows 257-264). Cut+Paste in one shot keeps
# the original per-cell styles (s="2"/"3"/"4") intact. Afterwards remove the
# now-empty leftover rows 3-9 on "New" (row 2 stays, ready for new data).
# ---------------------------------------------------------------------------
$destFirstRow = 257
$ws2.Range("A2:F9").Cut($ws1.Range("A" + $destFirstRow))
$ws2.Range("A3:F9").Delete(-4162)

# ---------------------------------------------------------------------------
# Step 3: re-create the hyperlinks on their new home (sheet "Previously
# added", column A, rows 257-264), then restore the original cell style
# that Hyperlinks.Add() overwrites (it forces the built-in "Hyperlink"
# style onto the cell).
# ---------------------------------------------------------------------------
$ws1.Range("A256").Copy() | Out-Null

for ($i = 0; $i -lt $oldLinks.Length; $i++) {
    $row = $destFirstRow + $i
    $cellRef = "A" + $row
    $ws1.Hyperlinks.Add($ws1.Range($cellRef), $oldLinks[$i])
    $ws1.Range($cellRef).PasteSpecial(-4122)
}

# ---------------------------------------------------------------------------
# Step 4: write the single new row of scraped data onto sheet "New", row 2.
# ---------------------------------------------------------------------------
$newUrl = "https://www.ss.com/msg/lv/real-estate/wood/balvi-and-reg/vecumu-pag/dkpmp.html"

$ws2.Range("A2").Value() = $newUrl
$ws2.Range("B2").Value() = "25 000 €"
$ws2.Range("C2").Value() = "Balvi un raj."
$ws2.Range("D2").Value() = "1.85 ha."
$ws2.Range("E2").Value() = ""
$ws2.Range("F2").Value() = 45978.76111111111

$ws2.Hyperlinks.Add($ws2.Range("A2"), $newUrl)

$ws1.Range("A256:F256").Copy() | Out-Null
$ws2.Range("A2:F2").PasteSpecial(-4122)

# Clean up the now-unused named "Hyperlink" cell style created by the
# Hyperlinks.Add() calls above (done once, at the very end, so that only a
# single extra font/style entry is ever created).
try {
    $wb.Styles.Item("Hyperlink").Delete()
} catch {
}
